$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.288.18'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.513.46'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.44'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.55%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.86'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.90%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.518.12'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.24%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.494'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.18%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.123'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.24%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.21'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.67%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.386'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.52%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.108.98'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.25%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000184'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.80%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.83'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -5.13%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.455.24'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.116'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '65.181.80'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.97'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -5.22%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.82'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.84%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.09'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '392.93'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.41%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.566'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.47%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('B24').Value = 'WrappedeETH'

$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.647.41'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.30'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.00%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('B26').Value = 'Dai'

$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('B27').Value = 'PEPE'

$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000112'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.42%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('B28').Value = 'RenderToken'

$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.58'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.84%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.996'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('B30').Value = 'PancakeSwap'

$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.23'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -7.44%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'

$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.18'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -7.31%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('B32').Value = 'RenzoRestakedETH'

$ws.Range('C32').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.528.21'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.38%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('B33').Value = 'USDe'

$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('B34').Value = 'EthereumClassic'

$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.30'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.07%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.11%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('B36').Value = 'Monero'

$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '171.09'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.69%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('B37').Value = 'Fetch.AI'

$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.22'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.92%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('B38').Value = 'Aptos'

$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.87'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.98%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('B39').Value = 'ImmutableX'

$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.51'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.21%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('B40').Value = 'NEARProtocol'

$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.79'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -6.28%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('B41').Value = 'Hedera'

$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0789'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.23%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('B42').Value = 'Mantle'

$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.816'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.27%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('B43').Value = 'EnergySwap'

$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.04'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +17.06%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.49'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.47%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('B46').Value = 'Filecoin'

$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.38'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -9.73%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('B47').Value = 'ONDO'

$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.19'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +6.91%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('B48').Value = 'Stacks'

$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.66'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.53%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('B49').Value = 'Cosmos'

$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.73'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.07%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('B50').Value = 'dogwifhat'

$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.10'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -8.23%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('B51').Value = 'Maker'

$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.243.31'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.30%  '
$ws.Range('E51').Style = 'Normal'
